$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: update with the newly scraped listing (keeps the row's slot) ---
$ws.Range("A2").Value = "2025-11-01 06:25:02"
$ws.Range("B2").Value = "【AWS→Xサーバー移行】サーバー構築の専門家募集"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5425003"
$ws.Range("G2").Value = 18
$ws.Range("H2").ClearContents()

# --- Row 3: update with the next newly scraped listing ---
$ws.Range("A3").Value = "2025-11-01 06:25:02"
$ws.Range("B3").Value = "【電子工作】基盤にDCケーブルのターミナルと抵抗を追加したい方募集!"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5424906"
$ws.Range("G3").Value = 10
$ws.Range("H3").ClearContents()

# --- Drop the old rows 4-10, which are no longer part of the refreshed list ---
$ws.Rows("4:10").Delete()

# --- Hyperlinks don't track cell/row edits automatically in this host, so
#     rebuild the collection from scratch: clear everything, then add back
#     links only for the two rows that remain. ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5425003")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5424906")

# Hyperlinks.Add stamps a freshly duplicated "Hyperlink" xf onto the cell;
# re-apply the named style so F2/F3 settle back on the original shared xf.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# --- Column B narrows slightly in the refreshed sheet. This host's
#     ColumnWidth setter round-trips through a pixel conversion that always
#     adds 5/6 back onto the stored OOXML width, so back that off here to
#     land exactly on width="36". ---
$ws.Columns("B").ColumnWidth = 36 - (5/6)
